$wb = $excel.ActiveWorkbook

# --- Sheet: Test Cases ---
$wsTest = $wb.Worksheets.Item("Test Cases")
$wsTest.Range("D2").Value2 = "PASS"
$wsTest.Range("D3").Value2 = "PASS"
$wsTest.Range("C5").Value2 = "N"
$wsTest.Range("D5").ClearContents()

# --- Sheet: PurchaseOrderCreation ---
$wsCreate = $wb.Worksheets.Item("PurchaseOrderCreation")
$wsCreate.Range("A2").Value2 = 144
$wsCreate.Range("B2").Value2 = "PI 144"

# --- Sheet: PurchaseOrderUpdate ---
$wsUpdate = $wb.Worksheets.Item("PurchaseOrderUpdate")
$wsUpdate.Range("A2").Value2 = 142
$wsUpdate.Range("B2").Value2 = "PI 142"

# --- Sheet: PurchaseOrderWorkflow ---
$wsFlow = $wb.Worksheets.Item("PurchaseOrderWorkflow")
$wsFlow.Range("A2").Value2 = 143
$wsFlow.Range("B2").Value2 = "PI 143"
$wsFlow.Range("C2").Value2 = 140
$wsFlow.Range("D2").Value2 = 1442
$wsFlow.Range("F2").Value2 = "PASS"

# --- Selections / active sheet (order matters: last Select() wins as active tab) ---
$wsTest.Range("B3").Select() | Out-Null
$wsUpdate.Range("C2").Select() | Out-Null
$wsFlow.Range("D2").Select() | Out-Null
